# Prefix each "command" name in column A with its protocol (sheet) name.
# Sheets LinaJourney, NRWaves, PersonalLina, PositiveSpin, ReEngagement are
# reference/config sheets and are left untouched. All other sheets get the
# sheet name prepended to every non-header value in column A.

$wb = $excel.ActiveWorkbook

$skipSheets = @("LinaJourney", "NRWaves", "PersonalLina", "PositiveSpin", "ReEngagement")

foreach ($ws in $wb.Worksheets) {
    if ($skipSheets -contains $ws.Name) {
        continue
    }

    $prefix = $ws.Name
    $usedRange = $ws.UsedRange
    $lastRow = $usedRange.Rows.Count

    for ($r = 2; $r -le $lastRow; $r++) {
        $cell = $ws.Cells.Item($r, 1)
        $current = $cell.Value2
        if ($current -ne $null -and $current -ne "") {
            $currentText = $current.ToString()
            if ($currentText.StartsWith("$prefix ") -eq $false) {
                $cell.Value2 = "$prefix $currentText"
            }
        }
    }
}
